# ---------------------------------------------------------------------------
# "Horarios actualizados Linea 141 - 400"
# Refresh the three scrape-snapshot sheets (LP1912, LP1912-215, 6203-6173)
# with the 18:49:07 pass: updates header metadata, fixes a handful of
# same-arrival-time tie-break swaps, and reshuffles/extends the tail of
# each sheet's "still approaching" bus queue (sorted by Hora_Llegada).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Write-RowBlock($ws, $startRow, $rowsList) {
    $nrows = $rowsList.Count
    $ncols = 5
    $arr = New-Object 'object[,]' $nrows,$ncols
    for ($r = 0; $r -lt $nrows; $r++) {
        for ($c = 0; $c -lt $ncols; $c++) {
            $arr[$r,$c] = $rowsList[$r][$c]
        }
    }
    $endRow = $startRow + $nrows - 1
    $target = $ws.Range($ws.Cells.Item($startRow,1), $ws.Cells.Item($endRow,5))
    $target.Value = $arr
}

# ===========================================================================
# Sheet 1: LP1912
# ===========================================================================
$ws = $wb.Worksheets.Item("LP1912")

$ws.Range("A2").Value = "Última actualización: 18:49:07"
$ws.Range("A3").Value = "Total filas: 426"

# Same-arrival-time (Hora_Llegada) row pairs whose tie-break order flipped
$ws.Cells.Item(71,3).Value = '215B_EL PATO'
$ws.Cells.Item(72,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(96,1).Value = '08:48:09'
$ws.Cells.Item(96,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(96,4).Value = 23
$ws.Cells.Item(97,1).Value = '08:36:20'
$ws.Cells.Item(97,3).Value = '16_SANTA ANA'
$ws.Cells.Item(97,4).Value = 35
$ws.Cells.Item(124,1).Value = '09:25:56'
$ws.Cells.Item(124,3).Value = '10_OLMOS'
$ws.Cells.Item(124,4).Value = 47
$ws.Cells.Item(125,1).Value = '08:19:33'
$ws.Cells.Item(125,3).Value = '15_ABASTO'
$ws.Cells.Item(125,4).Value = 113
$ws.Cells.Item(177,1).Value = '11:59:06'
$ws.Cells.Item(177,3).Value = '16_SANTA ANA'
$ws.Cells.Item(177,4).Value = 0
$ws.Cells.Item(178,1).Value = '10:12:35'
$ws.Cells.Item(178,3).Value = '225_GOMEZ'
$ws.Cells.Item(178,4).Value = 107
$ws.Cells.Item(190,1).Value = '10:52:48'
$ws.Cells.Item(190,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(190,4).Value = 77
$ws.Cells.Item(191,1).Value = '11:17:08'
$ws.Cells.Item(191,3).Value = '10_OLMOS'
$ws.Cells.Item(191,4).Value = 52
$ws.Cells.Item(202,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(203,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(284,1).Value = '14:44:25'
$ws.Cells.Item(284,3).Value = '17_ROMERO'
$ws.Cells.Item(284,4).Value = 72
$ws.Cells.Item(285,1).Value = '14:00:52'
$ws.Cells.Item(285,3).Value = '27_EL RETIRO'
$ws.Cells.Item(285,4).Value = 116
$ws.Cells.Item(325,3).Value = '16_SANTA ANA'
$ws.Cells.Item(326,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(365,1).Value = '17:42:01'
$ws.Cells.Item(365,3).Value = '15_ABASTO'
$ws.Cells.Item(365,4).Value = 34
$ws.Cells.Item(366,1).Value = '17:15:09'
$ws.Cells.Item(366,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(366,4).Value = 61

# Rows 384-431: tail of the queue re-sorted / extended with newly scraped
# buses (the 18:49:07 pass). Overwrite the full tail block at once.
$rows1 = New-Object System.Collections.Generic.List[object]
$rows1.Add(@('18:49:07','18:49','14X44_ABASTO',0,'LP1912')) | Out-Null
$rows1.Add(@('18:49:07','18:50','16_SANTA ANA',1,'LP1912')) | Out-Null
$rows1.Add(@('18:37:39','18:51','14_ABASTO',14,'LP1912')) | Out-Null
$rows1.Add(@('18:37:39','18:52','15_ABASTO',15,'LP1912')) | Out-Null
$rows1.Add(@('18:49:07','18:53','14_ABASTO',4,'LP1912')) | Out-Null
$rows1.Add(@('17:42:01','18:54','14_ABASTO',72,'LP1912')) | Out-Null
$rows1.Add(@('17:57:54','18:55','10_OLMOS',58,'LP1912')) | Out-Null
$rows1.Add(@('17:42:01','18:56','10_OLMOS',74,'LP1912')) | Out-Null
$rows1.Add(@('17:57:54','18:58','215A_EL PATO',61,'LP1912')) | Out-Null
$rows1.Add(@('17:15:09','18:59','215A_EL PATO',104,'LP1912')) | Out-Null
$rows1.Add(@('18:19:32','19:04','23_HERNANDEZ',45,'LP1912')) | Out-Null
$rows1.Add(@('17:57:54','19:04','11_ETCHEVERRY',67,'LP1912')) | Out-Null
$rows1.Add(@('17:15:09','19:05','11_ETCHEVERRY',110,'LP1912')) | Out-Null
$rows1.Add(@('18:37:39','19:05','23_HERNANDEZ',28,'LP1912')) | Out-Null
$rows1.Add(@('17:15:09','19:11','16_P MOR-SANTA ANA',116,'LP1912')) | Out-Null
$rows1.Add(@('17:42:01','19:12','10_OLMOS',90,'LP1912')) | Out-Null
$rows1.Add(@('18:49:07','19:16','15_ABASTO',27,'LP1912')) | Out-Null
$rows1.Add(@('17:57:54','19:16','27_EL RETIRO',79,'LP1912')) | Out-Null
$rows1.Add(@('17:57:54','19:17','27_EL RETIRO',95,'LP1912')) | Out-Null
$rows1.Add(@('17:42:01','19:17','16_P MOR-SANTA ANA',95,'LP1912')) | Out-Null
$rows1.Add(@('17:57:54','19:20','26_HERNANDEZ',83,'LP1912')) | Out-Null
$rows1.Add(@('17:57:54','19:20','14_ABASTO',83,'LP1912')) | Out-Null
$rows1.Add(@('17:42:01','19:21','26_HERNANDEZ',99,'LP1912')) | Out-Null
$rows1.Add(@('18:49:07','19:21','16_SANTA ANA',32,'LP1912')) | Out-Null
$rows1.Add(@('18:49:07','19:21','14_ABASTO',32,'LP1912')) | Out-Null
$rows1.Add(@('17:57:54','19:29','225_GOMEZ',92,'LP1912')) | Out-Null
$rows1.Add(@('17:42:01','19:30','225_GOMEZ',108,'LP1912')) | Out-Null
$rows1.Add(@('17:57:54','19:31','16_P MOR-SANTA ANA',94,'LP1912')) | Out-Null
$rows1.Add(@('18:37:39','19:34','23_HERNANDEZ',57,'LP1912')) | Out-Null
$rows1.Add(@('18:49:07','19:35','23_HERNANDEZ',46,'LP1912')) | Out-Null
$rows1.Add(@('17:57:54','19:39','215C_EL PATO',102,'LP1912')) | Out-Null
$rows1.Add(@('17:42:01','19:40','215C_EL PATO',118,'LP1912')) | Out-Null
$rows1.Add(@('17:57:54','19:49','11X44_ETCHEVERRY',112,'LP1912')) | Out-Null
$rows1.Add(@('17:57:54','19:50','16_P MOR-SANTA ANA',113,'LP1912')) | Out-Null
$rows1.Add(@('18:19:32','19:50','11X44_ETCHEVERRY',91,'LP1912')) | Out-Null
$rows1.Add(@('18:49:07','19:51','81_EL PELIGRO',62,'LP1912')) | Out-Null
$rows1.Add(@('18:49:07','19:51','16_P MOR-SANTA ANA',62,'LP1912')) | Out-Null
$rows1.Add(@('18:19:32','19:54','16_P MOR-SANTA ANA',95,'LP1912')) | Out-Null
$rows1.Add(@('18:19:32','19:59','17_ROMERO',100,'LP1912')) | Out-Null
$rows1.Add(@('18:19:32','20:10','16_P MOR-167 Y 521',111,'LP1912')) | Out-Null
$rows1.Add(@('18:49:07','20:11','16_P MOR-167 Y 521',82,'LP1912')) | Out-Null
$rows1.Add(@('18:37:39','20:12','16_P MOR-SANTA ANA',95,'LP1912')) | Out-Null
$rows1.Add(@('18:37:39','20:20','26_HERNANDEZ',103,'LP1912')) | Out-Null
$rows1.Add(@('18:49:07','20:21','26_HERNANDEZ',92,'LP1912')) | Out-Null
$rows1.Add(@('18:37:39','20:22','11_ETCHEVERRY',105,'LP1912')) | Out-Null
$rows1.Add(@('18:37:39','20:23','215A_EL PATO',106,'LP1912')) | Out-Null
$rows1.Add(@('18:49:07','20:24','215A_EL PATO',95,'LP1912')) | Out-Null
$rows1.Add(@('18:37:39','20:31','225_GOMEZ',114,'LP1912')) | Out-Null
Write-RowBlock $ws 384 $rows1

# ===========================================================================
# Sheet 2: LP1912-215
# ===========================================================================
$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Range("A2").Value = "Última actualización: 18:49:07"
$ws.Range("A3").Value = "Total filas: 42"

$rows2 = New-Object System.Collections.Generic.List[object]
$rows2.Add(@('18:49:07','20:24','215A_EL PATO',95,'LP1912')) | Out-Null
Write-RowBlock $ws 47 $rows2

# ===========================================================================
# Sheet 3: 6203-6173
# ===========================================================================
$ws = $wb.Worksheets.Item("6203-6173")

$ws.Range("A2").Value = "Última actualización: 18:49:07"
$ws.Range("A3").Value = "Total filas: 55"

# A new bus is inserted at row 53 (pushing the former rows 53-58 down to
# 54-59) and another new bus is appended at the end as row 60.
$rows3 = New-Object System.Collections.Generic.List[object]
$rows3.Add(@('18:49:07','18:50','215A_LA PLATA',1,'L6173')) | Out-Null
$rows3.Add(@('17:57:54','18:51','215A_LA PLATA',54,'L6173')) | Out-Null
$rows3.Add(@('17:15:09','18:52','215A_LA PLATA',97,'L6173')) | Out-Null
$rows3.Add(@('17:57:54','19:03','215B_LP-P MOR-1 Y 57',66,'L6173')) | Out-Null
$rows3.Add(@('17:15:09','19:04','215B_LP-P MOR-1 Y 57',109,'L6173')) | Out-Null
$rows3.Add(@('18:37:39','19:53','215C_LA PLATA',76,'L6203')) | Out-Null
$rows3.Add(@('18:19:32','19:54','215C_LA PLATA',95,'L6203')) | Out-Null
$rows3.Add(@('18:49:07','20:39','215A_LA PLATA',110,'L6173')) | Out-Null
Write-RowBlock $ws 53 $rows3
